$wb = $excel.ActiveWorkbook
Write-Output ($wb | Get-Member)
